$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells first so that numeric-looking
# strings (e.g. "4.80", "0.480") are preserved exactly as text, matching
# the source data which stores these values as literal text, not numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

$ws.Range('D2').Value = '68.405.12'
$ws.Range('E2').Value = '  +1.03%  '
$ws.Range('D3').Value = '2.436.96'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '557.72'
$ws.Range('E5').Value = '  +1.64%  '
$ws.Range('E6').Value = '  +2.34%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +3.22%  '
$ws.Range('E9').Value = '  +9.53%  '
$ws.Range('E10').Value = '  +0.46%  '
$ws.Range('D11').Value = '4.80'
$ws.Range('E11').Value = '  +2.19%  '
$ws.Range('D12').Value = '0.327'
$ws.Range('E12').Value = '  -1.41%  '
$ws.Range('D13').Value = '68.306.02'
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('E14').Value = '  +4.18%  '
$ws.Range('D15').Value = '23.24'
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('D16').Value = '10.41'
$ws.Range('E16').Value = '  -2.21%  '
$ws.Range('D17').Value = '336.89'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').Value = '6.89'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('E19').Value = '  +2.33%  '
$ws.Range('E20').Value = '  +3.98%  '
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').Value = '66.75'
$ws.Range('E22').Value = '  +1.33%  '
$ws.Range('D23').Value = '3.68'
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('D24').Value = '8.14'
$ws.Range('E24').Value = '  +2.63%  '
$ws.Range('D25').Value = '0.0₃0814'
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('D26').Value = '7.18'
$ws.Range('E26').Value = '  +2.96%  '
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('D28').Value = '425.05'
$ws.Range('E28').Value = '  +1.42%  '
$ws.Range('E29').Value = '  +2.83%  '
$ws.Range('E30').Value = '  +1.10%  '
$ws.Range('D31').Value = '160.55'
$ws.Range('E31').Value = '  +2.58%  '
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').Value = '17.77'
$ws.Range('E34').Value = '  +1.62%  '
$ws.Range('E35').Value = '  -1.55%  '
$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D36').Value = '0.296'
$ws.Range('E36').Value = '  -0.29%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '4.36'
$ws.Range('E37').Value = '  +2.28%  '
$ws.Range('E38').Value = '  +3.48%  '
$ws.Range('D39').Value = '1.06'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').Value = '2.01'
$ws.Range('E40').Value = '  +1.63%  '
$ws.Range('D41').Value = '3.35'
$ws.Range('E41').Value = '  +2.81%  '
$ws.Range('D42').Value = '129.46'
$ws.Range('E42').Value = '  -1.86%  '
$ws.Range('D43').Value = '0.0716'
$ws.Range('E43').Value = '  +1.24%  '
$ws.Range('D44').Value = '0.480'
$ws.Range('E44').Value = '  +2.32%  '
$ws.Range('D45').Value = '0.561'
$ws.Range('E46').Value = '  +2.29%  '
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('D49').Value = '4.90'
$ws.Range('E49').Value = '  -2.89%  '
$ws.Range('D50').Value = '16.68'
$ws.Range('E50').Value = '  +1.47%  '
$ws.Range('E51').Value = '  +5.48%  '
